$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Header summary updates
# ------------------------------------------------------------------
$ws.Range("E11").Value = 382042       # VALOR MORA total
$ws.Range("C13").Value = 2            # Cant. Trabajadores
$ws.Range("F13").Value = 8            # Cant. Periodos

# ------------------------------------------------------------------
# 2. Insert five blank rows right before the old "last" worker row
#    (row 22). This naturally pushes that row down to 27 (keeping
#    its distinctive bottom-border styling) and pushes the footer
#    block down from rows 27/28 to 32/33, updating merged ranges
#    automatically.
# ------------------------------------------------------------------
$ws.Rows("22:26").Insert()

# Stamp the plain "middle row" styling (as used by row 16) onto the
# newly inserted rows 22-26.
$ws.Range("B16:J16").Copy() | Out-Null
$ws.Range("B22:J26").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3. Re-populate the worker data block (rows 16-27) in strict row
#    order so the shared-string table is rebuilt in the same order
#    the source workbook uses. First clear the old period values so
#    they are dropped from the shared-string table, then re-add them
#    in the new (descending) order.
# ------------------------------------------------------------------
$ws.Range("E16:E21").ClearContents()
$ws.Cells.Item(27, 5).ClearContents()

$ws.Cells.Item(16, 5).Value = "1911"
$ws.Cells.Item(16, 7).Value = 828116

$ws.Cells.Item(17, 5).Value = "1910"
$ws.Cells.Item(17, 7).Value = 828116

$ws.Cells.Item(18, 5).Value = "1909"
$ws.Cells.Item(18, 7).Value = 828116

$ws.Cells.Item(19, 5).Value = "1908"
$ws.Cells.Item(19, 7).Value = 828116

$ws.Cells.Item(20, 5).Value = "1907"
$ws.Cells.Item(20, 7).Value = 828116

$ws.Cells.Item(21, 5).Value = "1906"
$ws.Cells.Item(21, 7).Value = 828116

$ws.Cells.Item(22, 2).Value = "CC"
$ws.Cells.Item(22, 3).Value = "73130670"
$ws.Cells.Item(22, 4).Value = "ROBERTO ROZO ACUNA"
$ws.Cells.Item(22, 5).Value = "1905"
$ws.Cells.Item(22, 6).Value = 33125
$ws.Cells.Item(22, 7).Value = 828116

$ws.Cells.Item(23, 2).Value = "CC"
$ws.Cells.Item(23, 3).Value = "1019065959"
$ws.Cells.Item(23, 4).Value = "GUIDO ANDRES TORRES RODRIGUEZ"
$ws.Cells.Item(23, 5).Value = "1912"
$ws.Cells.Item(23, 6).Value = 17667
$ws.Cells.Item(23, 7).Value = 828116

$ws.Cells.Item(24, 2).Value = "CC"
$ws.Cells.Item(24, 3).Value = "1019065959"
$ws.Cells.Item(24, 4).Value = "GUIDO ANDRES TORRES RODRIGUEZ"
$ws.Cells.Item(24, 5).Value = "1911"
$ws.Cells.Item(24, 6).Value = 33125
$ws.Cells.Item(24, 7).Value = 828116

$ws.Cells.Item(25, 2).Value = "CC"
$ws.Cells.Item(25, 3).Value = "1019065959"
$ws.Cells.Item(25, 4).Value = "GUIDO ANDRES TORRES RODRIGUEZ"
$ws.Cells.Item(25, 5).Value = "1910"
$ws.Cells.Item(25, 6).Value = 33125
$ws.Cells.Item(25, 7).Value = 828116

$ws.Cells.Item(26, 2).Value = "CC"
$ws.Cells.Item(26, 3).Value = "1019065959"
$ws.Cells.Item(26, 4).Value = "GUIDO ANDRES TORRES RODRIGUEZ"
$ws.Cells.Item(26, 5).Value = "1909"
$ws.Cells.Item(26, 6).Value = 33125
$ws.Cells.Item(26, 7).Value = 828116

$ws.Cells.Item(27, 3).Value = "1019065959"
$ws.Cells.Item(27, 4).Value = "GUIDO ANDRES TORRES RODRIGUEZ"
$ws.Cells.Item(27, 5).Value = "1908"
$ws.Cells.Item(27, 6).Value = 33125
$ws.Cells.Item(27, 7).Value = 828116

# ------------------------------------------------------------------
# 4. Widen column D to fit the longer worker name.
# ------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 34.1796875
